$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update D (Price) and E (Volume 1h) columns with refreshed crypto data ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.762.85"
$ws.Range("E2").Value = "  +0.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.817.05"
$ws.Range("E3").Value = "  +1.77%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.24"
$ws.Range("E5").Value = "  +1.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.01"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("E9").Value = "  +0.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.30"
$ws.Range("E10").Value = "  -0.51%  "

$ws.Range("E11").Value = "  +0.85%  "

$ws.Range("E12").Value = "  -0.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.76"
$ws.Range("E13").Value = "  -0.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.456.49"
$ws.Range("E14").Value = "  +1.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.810.17"
$ws.Range("E15").Value = "  +1.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.782.33"
$ws.Range("E16").Value = "  +0.68%  "

$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("E18").Value = "  +1.27%  "

$ws.Range("E19").Value = "  +0.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "462.87"
$ws.Range("E20").Value = "  +1.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.89"
$ws.Range("E21").Value = "  -0.97%  "

$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000148"
$ws.Range("E23").Value = "  -3.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.36"

$ws.Range("E25").Value = "  +1.74%  "

$ws.Range("E26").Value = "  -1.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.04"
$ws.Range("E27").Value = "  -0.73%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.966.66"
$ws.Range("E29").Value = "  +1.77%  "

$ws.Range("E30").Value = "  +0.67%  "

$ws.Range("E31").Value = "  +1.63%  "

$ws.Range("E32").Value = "  +2.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.56"

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("E35").Value = "  -0.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0999"
$ws.Range("E36").Value = "  -0.21%  "

$ws.Range("E37").Value = "  -0.81%  "

$ws.Range("E38").Value = "  +0.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("E40").Value = "  +1.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.06"
$ws.Range("E43").Value = "  +2.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.14"
$ws.Range("E46").Value = "  -4.51%  "

$ws.Range("E47").Value = "  +12.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "148.82"
$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("E50").Value = "  +0.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "386.64"
$ws.Range("E51").Value = "  -0.66%  "

# --- Rows 44 and 45 swapped coin entries (EnergySwap <-> TheGraph) with new data ---
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.300"
$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "28.41"
$ws.Range("E45").Value = "  +8.95%  "
